$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "147×6=" "582×5="
Replace-Text "213×5=" "705×5="
Replace-Text "130×8=" "440×6="
Replace-Text "209×7=" "141×8="
Replace-Text "250×7=" "384×4="
Replace-Text "648×8=" "965×3="
Replace-Text "221×4=" "732×3="
Replace-Text "910×5=" "856×3="
Replace-Text "761×4=" "694×7="
Replace-Text "743×6=" "209×3="
Replace-Text "887×5=" "612×4="
Replace-Text "983×2=" "875×5="
Replace-Text "341×7=" "625×5="
Replace-Text "111×9=" "261×5="
Replace-Text "922×7=" "246×2="
Replace-Text "911×2=" "353×2="
Replace-Text "908×2=" "313×9="
Replace-Text "169×5=" "574×3="
Replace-Text "919×7=" "452×5="
Replace-Text "311×2=" "297×7="
Replace-Text "249×7=" "796×4="
Replace-Text "864×5=" "927×5="
Replace-Text "751×6=" "364×8="
Replace-Text "684×8=" "991×9="
Replace-Text "114×2=" "260×8="
